$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 304-305 (existing rows 304.. shift down to 306..)
$ws.Range("A304:A305").EntireRow.Insert()

# New row 304: Conconina(o), Región Metropolitana, week of 2021-09-10
$ws.Range("A304").Value = 11
$ws.Range("B304").Value = "Vega Monumental Concepción"
$ws.Range("C304").Value = "Bíobío"
$ws.Range("D304").Value = 44449
$ws.Range("E304").Value = 8
$ws.Range("F304").Value = 100112033
$ws.Range("G304").Value = "Lechuga"
$ws.Range("H304").Value = "Conconina(o)"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 100
$ws.Range("K304").Value = 6000
$ws.Range("L304").Value = 6500
$ws.Range("M304").Value = 6250
$ws.Range("N304").Value = "$/caja 10 unidades"
$ws.Range("O304").Value = "Región Metropolitana"
$ws.Range("P304").Value = 625
$ws.Range("Q304").Value = 10
$ws.Range("R304").Value = "Hortaliza"

# New row 305: Escarola, Región de Coquimbo, week of 2021-09-10
$ws.Range("A305").Value = 11
$ws.Range("B305").Value = "Vega Monumental Concepción"
$ws.Range("C305").Value = "Bíobío"
$ws.Range("D305").Value = 44449
$ws.Range("E305").Value = 8
$ws.Range("F305").Value = 100112033
$ws.Range("G305").Value = "Lechuga"
$ws.Range("H305").Value = "Escarola"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 100
$ws.Range("K305").Value = 8000
$ws.Range("L305").Value = 8500
$ws.Range("M305").Value = 8250
$ws.Range("N305").Value = "$/caja 15 unidades"
$ws.Range("O305").Value = "Región de Coquimbo"
$ws.Range("P305").Value = 550
$ws.Range("Q305").Value = 15
$ws.Range("R305").Value = "Hortaliza"
